$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Malasia: moves from after "Polonia" to right after "Chile" (row 30), with updated stats ---
$ws.Rows("30:30").Insert()
$ws.Range("A30").Value = "Malasia"
$ws.Range("B30").Value = 3333
$ws.Range("C30").Value = 217
$ws.Range("D30").Value = 827
$ws.Range("E30").Value = 2453
$ws.Range("F30").Value = 105
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 53
$ws.Rows("33:33").Delete()

# --- Libano: moves from after "Moldavia" to right after "Bosnia y Herzegovina" (row 72), with updated stats ---
$ws.Rows("72:72").Insert()
$ws.Range("A72").Value = "Libano"
$ws.Range("B72").Value = 508
$ws.Range("C72").Value = 14
$ws.Range("D72").Value = 46
$ws.Range("E72").Value = 445
$ws.Range("F72").Value = 3
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 17
$ws.Rows("74:74").Delete()

# --- Isla de Man: moves from after "Trinidad yTobago" to right after "Mayotte" (row 119), with updated stats ---
$ws.Rows("119:119").Insert()
$ws.Range("A119").Value = "Isla de Man"
$ws.Range("B119").Value = 114
$ws.Range("C119").Value = 19
$ws.Range("D119").Value = 0
$ws.Range("E119").Value = 113
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 1
$ws.Rows("124:124").Delete()

# --- Plain numeric updates (no row movement) ---
# Suiza (row 12)
$ws.Range("B12").Value = 19106
$ws.Range("C12").Value = 279
$ws.Range("D12").Value = 4846
$ws.Range("E12").Value = 13695
$ws.Range("G12").Value = 29
$ws.Range("H12").Value = 565

# Belgica (row 14)
$ws.Range("B14").Value = 16770
$ws.Range("C14").Value = 1422
$ws.Range("D14").Value = 2872
$ws.Range("E14").Value = 12755
$ws.Range("F14").Value = 1205
$ws.Range("G14").Value = 132
$ws.Range("H14").Value = 1143

# Austria (row 17)
$ws.Range("B17").Value = 11238
$ws.Range("C17").Value = 109
$ws.Range("E17").Value = 9058

# Singapur (row 53)
$ws.Range("B53").Value = 1114
$ws.Range("C53").Value = 65
$ws.Range("E53").Value = 843

# Hong Kong (row 63)
$ws.Range("B63").Value = 845
$ws.Range("C63").Value = 43
$ws.Range("E63").Value = 687

# --- Timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 11:20"
